$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "85-43="
$t.Cell(1,2).Range.Text = "44+36="
$t.Cell(1,3).Range.Text = "48+27="
$t.Cell(1,4).Range.Text = "18+77="
$t.Cell(1,5).Range.Text = "33+29="
$t.Cell(2,1).Range.Text = "25+59="
$t.Cell(2,2).Range.Text = "19+1="
$t.Cell(2,3).Range.Text = "32-14="
$t.Cell(2,4).Range.Text = "55+41="
$t.Cell(2,5).Range.Text = "79-56="
$t.Cell(3,1).Range.Text = "60-14="
$t.Cell(3,2).Range.Text = "46+10="
$t.Cell(3,3).Range.Text = "37+28="
$t.Cell(3,4).Range.Text = "70-46="
$t.Cell(3,5).Range.Text = "9+61="
$t.Cell(4,1).Range.Text = "28+35="
$t.Cell(4,2).Range.Text = "35+45="
$t.Cell(4,3).Range.Text = "65-47="
$t.Cell(4,4).Range.Text = "95-21="
$t.Cell(4,5).Range.Text = "15-3="
$t.Cell(5,1).Range.Text = "59+9="
$t.Cell(5,2).Range.Text = "54-0="
$t.Cell(5,3).Range.Text = "83-4="
$t.Cell(5,4).Range.Text = "95-23="
$t.Cell(5,5).Range.Text = "63-12="
$t.Cell(6,1).Range.Text = "37+30="
$t.Cell(6,2).Range.Text = "14+53="
$t.Cell(6,3).Range.Text = "83-6="
$t.Cell(6,4).Range.Text = "13+52="
$t.Cell(6,5).Range.Text = "40-16="
$t.Cell(7,1).Range.Text = "13+25="
$t.Cell(7,2).Range.Text = "87-60="
$t.Cell(7,3).Range.Text = "50-24="
$t.Cell(7,4).Range.Text = "84-70="
$t.Cell(7,5).Range.Text = "27+55="
$t.Cell(8,1).Range.Text = "97-88="
$t.Cell(8,2).Range.Text = "35+19="
$t.Cell(8,3).Range.Text = "19+30="
$t.Cell(8,4).Range.Text = "28-19="
$t.Cell(8,5).Range.Text = "91-79="
$t.Cell(9,1).Range.Text = "14+23="
$t.Cell(9,2).Range.Text = "64+17="
$t.Cell(9,3).Range.Text = "57-33="
$t.Cell(9,4).Range.Text = "22-13="
$t.Cell(9,5).Range.Text = "56+19="
$t.Cell(10,1).Range.Text = "9+81="
$t.Cell(10,2).Range.Text = "4+38="
$t.Cell(10,3).Range.Text = "30+12="
$t.Cell(10,4).Range.Text = "28+45="
$t.Cell(10,5).Range.Text = "1+48="
$t.Cell(11,1).Range.Text = "97-7="
$t.Cell(11,2).Range.Text = "92-1="
$t.Cell(11,3).Range.Text = "54+23="
$t.Cell(11,4).Range.Text = "85-59="
$t.Cell(11,5).Range.Text = "88-55="
$t.Cell(12,1).Range.Text = "75-14="
$t.Cell(12,2).Range.Text = "70+3="
$t.Cell(12,3).Range.Text = "88-48="
$t.Cell(12,4).Range.Text = "92-60="
$t.Cell(12,5).Range.Text = "11+64="
$t.Cell(13,1).Range.Text = "91-35="
$t.Cell(13,2).Range.Text = "38-17="
$t.Cell(13,3).Range.Text = "95-26="
$t.Cell(13,4).Range.Text = "50+12="
$t.Cell(13,5).Range.Text = "96-64="
$t.Cell(14,1).Range.Text = "69-49="
$t.Cell(14,2).Range.Text = "0+47="
$t.Cell(14,3).Range.Text = "89-15="
$t.Cell(14,4).Range.Text = "33-11="
$t.Cell(14,5).Range.Text = "28+41="
$t.Cell(15,1).Range.Text = "77-29="
$t.Cell(15,2).Range.Text = "27+36="
$t.Cell(15,3).Range.Text = "39+50="
$t.Cell(15,4).Range.Text = "60-11="
$t.Cell(15,5).Range.Text = "22-10="
$t.Cell(16,1).Range.Text = "11+62="
$t.Cell(16,2).Range.Text = "27+3="
$t.Cell(16,3).Range.Text = "12+54="
$t.Cell(16,4).Range.Text = "41+44="
$t.Cell(16,5).Range.Text = "89-68="
$t.Cell(17,1).Range.Text = "77+20="
$t.Cell(17,2).Range.Text = "40+56="
$t.Cell(17,3).Range.Text = "23-13="
$t.Cell(17,4).Range.Text = "83-35="
$t.Cell(17,5).Range.Text = "48-1="
$t.Cell(18,1).Range.Text = "88-58="
$t.Cell(18,2).Range.Text = "26+6="
$t.Cell(18,3).Range.Text = "98-91="
$t.Cell(18,4).Range.Text = "68+19="
$t.Cell(18,5).Range.Text = "56-8="
$t.Cell(19,1).Range.Text = "23-10="
$t.Cell(19,2).Range.Text = "83-38="
$t.Cell(19,3).Range.Text = "86-83="
$t.Cell(19,4).Range.Text = "62+17="
$t.Cell(19,5).Range.Text = "43+45="
$t.Cell(20,1).Range.Text = "26+24="
$t.Cell(20,2).Range.Text = "70+14="
$t.Cell(20,3).Range.Text = "49+5="
$t.Cell(20,4).Range.Text = "63-32="
$t.Cell(20,5).Range.Text = "30+2="
